# Generate Report for Handoff
# Adds a new localization-status row for file
# "f9b55918-926f-4669-8b21-1c9a060826f3.md" to the Overview, zh-cn and
# de-de sheets (mirroring the existing 472ab7fd-... row), expanding each
# sheet's table by one row.

$wb = $excel.ActiveWorkbook

$fileGuid = "f9b55918-926f-4669-8b21-1c9a060826f3"
$commitSha = "b691d73476a77a285f97898ad2f30cc86bf1e7bc"
$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68ce46539d6a7e3e75e2c06b44d47a1335bfcb06/e2e/"

# ---------------------------------------------------------------------
# Overview sheet: File Name | Path And Name | Extension | Publish URL |
#                 zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "$fileGuid.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-05 06:45:26"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "$repoBase$fileGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\$fileGuid.md"
) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Source Path |
#              Priority | Content Duplicate | Latest Handoff File |
#              Latest Handoff Datetime | Latest Target File |
#              Latest Handback File | Latest Handback DateTime |
#              Reference Tokens | To be localized | Dependency From |
#              Has metadata | Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "$fileGuid.$commitSha.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 06:45:21"
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("J3").Value = "'"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "$repoBase$fileGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$fileGuid.md"
) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: same column layout as zh-cn
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "$fileGuid.$commitSha.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 06:45:26"
$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("J3").Value = "'"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "$repoBase$fileGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$fileGuid.md"
) | Out-Null
